# Updated the workbook so the login-test sheet reads its credentials from
# Excel data: add a new "ValidLogin" sheet (UserName/Password header row +
# admin/manager data row) and remove the old placeholder "Sheet1".

# Avoid the "are you sure you want to delete this sheet" prompt when the
# original sheet (which still has data on it) is removed below.
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Add the replacement sheet and populate it with the login fixture data.
$validLogin = $wb.Worksheets.Add()
$validLogin.Name = "ValidLogin"

$validLogin.Range("A1").Value = "UserName"
$validLogin.Range("B1").Value = "Password"
$validLogin.Range("A2").Value = "admin"
$validLogin.Range("B2").Value = "manager"

# Drop the old "Sheet1" now that "ValidLogin" carries the data.
[void]$wb.Worksheets.Item("Sheet1").Delete()

# Match the saved view state: selection on B3, zoomed to 175%.
[void]$validLogin.Range("B3").Select()
$excel.ActiveWindow.Zoom = 175

$excel.DisplayAlerts = $true
